# Add task estimation for the "ChangePassword" user story to the
# Estimacija workbook. This appends one user-story row (col A) followed
# by four task rows (col B = task, col C = estimated time in minutes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user story header (row 52, column A only - matches style of other
# user-story rows such as A48).
$ws.Range("A52").Value = "Kao korisnik potrebno je da mogu da promenim svoj PASS"

# New tasks under the story (rows 53-56): task description in column B,
# estimated time (minutes) in column C. No "real time" yet since the
# tasks are only estimates.
$ws.Range("B53").Value = "Dodavanje ChangePassword metode u Wallet I WalletService"
$ws.Range("C53").Value = 10

$ws.Range("B54").Value = "Implementacija testova za ChangePassword"
$ws.Range("C54").Value = 15

$ws.Range("B55").Value = "Dodavanje rute za ChangePassword na WalletController"
$ws.Range("C55").Value = 5

$ws.Range("B56").Value = "Dodavanje stranice za ChangePassword u MVC aplikaciju"
$ws.Range("C56").Value = 20

# Move the active selection to the next empty cell in column D, matching
# where the author left off after entering the new estimates.
$ws.Range("D53").Select()
